# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Mateus_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as described by the commit 'chore: update Sheets via scheduled runner'.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")

# Row 64: Forged from the Void
$ws.Range("H64").Value = 8091.4614
$ws.Range("I64").Value = 4129.6665
$ws.Range("J64").Value = 9280
$ws.Range("K64").Value = 4129.6665
$ws.Range("L64").Value = 9280
$ws.Range("M64").Value = -3881.6665
$ws.Range("N64").Value = -9776

# Row 67: Dodging the Draft (L)
$ws.Range("H67").Value = 8091.4614
$ws.Range("I67").Value = 4129.6665
$ws.Range("J67").Value = 9280
$ws.Range("K67").Value = 4129.6665
$ws.Range("L67").Value = 9280
$ws.Range("M67").Value = -3271.6665
$ws.Range("N67").Value = -10996

# Row 87: There Was a Late Fee
$ws.Range("H87").Value = 57400
$ws.Range("J87").Value = 57400
$ws.Range("L87").Value = 57400
$ws.Range("N87").Value = -59896

# Row 90: A Gate Arcane Is Dragon's Bane (L)
$ws.Range("H90").Value = 57400
$ws.Range("J90").Value = 57400
$ws.Range("L90").Value = 172200
$ws.Range("N90").Value = -184680

# Row 116: Growing Up
$ws.Range("H116").Value = 4033
$ws.Range("J116").Value = 4900
$ws.Range("L116").Value = 4900
$ws.Range("N116").Value = -11784

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 913.6923
$ws.Range("I132").Value = 960.125
$ws.Range("K132").Value = 2880.375
$ws.Range("M132").Value = -350.375

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")

# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 2566.6667
$ws.Range("I2").Value = 2566.6667
$ws.Range("K2").Value = 2566.6667
$ws.Range("M2").Value = -2453.6667

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 9449.5
$ws.Range("I45").Value = 7500
$ws.Range("K45").Value = 7500
$ws.Range("M45").Value = -7123

# Row 55: Employee Retention
$ws.Range("H55").Value = 14800
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 6679.5293
$ws.Range("I110").Value = 4879.5
$ws.Range("J110").Value = 10999.6
$ws.Range("K110").Value = 4879.5
$ws.Range("L110").Value = 10999.6
$ws.Range("M110").Value = -2834.5
$ws.Range("N110").Value = -15089.6

# Row 116: No Scope
$ws.Range("H116").Value = 2566.6667
$ws.Range("I116").Value = 2566.6667
$ws.Range("K116").Value = 2566.6667
$ws.Range("M116").Value = -272.6667000000002

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 4809.3335
$ws.Range("J132").Value = 15000
$ws.Range("L132").Value = 45000
$ws.Range("N132").Value = -50060

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")

# Row 3: Hells Bells
$ws.Range("H3").Value = 2566.6667
$ws.Range("I3").Value = 2566.6667
$ws.Range("K3").Value = 2566.6667
$ws.Range("M3").Value = -2452.6667

# Row 22: Riveting Run
$ws.Range("H22").Value = 4074
$ws.Range("I22").Value = 5000
$ws.Range("K22").Value = 5000
$ws.Range("M22").Value = -4827

# Row 82: Spirituality Inspector
$ws.Range("H82").Value = 20051.4
$ws.Range("J82").Value = 42000
$ws.Range("L82").Value = 42000
$ws.Range("N82").Value = -42766

# Row 85: The Clamor for Hammers (L)
$ws.Range("H85").Value = 20051.4
$ws.Range("J85").Value = 42000
$ws.Range("L85").Value = 42000
$ws.Range("N85").Value = -44652

# Row 129: Pruned to Perfection
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")

# Row 41: The Lone Bowman
$ws.Range("H41").Value = 16812.572
$ws.Range("I41").Value = 15677.6
$ws.Range("J41").Value = 19650
$ws.Range("K41").Value = 15677.6
$ws.Range("L41").Value = 19650
$ws.Range("M41").Value = -15249.6
$ws.Range("N41").Value = -20506

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 4384.222
$ws.Range("I58").Value = 1732.6957
$ws.Range("J58").Value = 9075.385
$ws.Range("K58").Value = 1732.6957
$ws.Range("L58").Value = 9075.385
$ws.Range("M58").Value = -1529.6957
$ws.Range("N58").Value = -9481.385

# Row 59: Bow Down to Magic
$ws.Range("H59").Value = 41100
$ws.Range("J59").Value = 41100
$ws.Range("L59").Value = 41100
$ws.Range("N59").Value = -43390

# Row 62: Splinter in the Sewers
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# Row 65: The Lumber of Their Discontent (L)
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# Row 74: License to Heal
$ws.Range("H74").Value = 27334.889
$ws.Range("I74").Value = 19980
$ws.Range("J74").Value = 28254.25
$ws.Range("K74").Value = 19980
$ws.Range("L74").Value = 28254.25
$ws.Range("M74").Value = -19106
$ws.Range("N74").Value = -30002.25

# Row 77: Purified Polyrhythm (L)
$ws.Range("H77").Value = 27334.889
$ws.Range("I77").Value = 19980
$ws.Range("J77").Value = 28254.25
$ws.Range("K77").Value = 59940
$ws.Range("L77").Value = 84762.75
$ws.Range("M77").Value = -55572
$ws.Range("N77").Value = -93498.75

# Row 97: Wood That You Could
$ws.Range("H97").Value = 33333.6
$ws.Range("J97").Value = 34167
$ws.Range("L97").Value = 34167
$ws.Range("N97").Value = -36149

# Row 136: Turali Quality
$ws.Range("H136").Value = 4384.222
$ws.Range("I136").Value = 1732.6957
$ws.Range("J136").Value = 9075.385
$ws.Range("K136").Value = 5198.0871
$ws.Range("L136").Value = 27226.155
$ws.Range("M136").Value = -2648.0871
$ws.Range("N136").Value = -32326.155

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")

# Row 3: Trout Fishing in Limsa
$ws.Range("H3").Value = 4825.222
$ws.Range("I3").Value = 4947.0713
$ws.Range("K3").Value = 14841.2139
$ws.Range("M3").Value = -14729.2139

# Row 8: Whip It
$ws.Range("H8").Value = 1256
$ws.Range("I8").Value = 1256
$ws.Range("K8").Value = 3768
$ws.Range("M8").Value = -3629

# Row 11: Putting the Squeeze On
$ws.Range("H11").Value = 20000440
$ws.Range("I11").Value = 20000440
$ws.Range("K11").Value = 60001320
$ws.Range("M11").Value = -60001180

# Row 122: Salt of the North
$ws.Range("H122").Value = 1025.3334
$ws.Range("J122").Value = 1438.5
$ws.Range("L122").Value = 12946.5
$ws.Range("N122").Value = -17846.5

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 7145357
$ws.Range("J131").Value = 7694922.5
$ws.Range("L131").Value = 23084767.5
$ws.Range("N131").Value = -23094847.5

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")

# Row 34: All Booked Up
$ws.Range("H34").Value = 69000
$ws.Range("J34").Value = 69000
$ws.Range("L34").Value = 69000
$ws.Range("N34").Value = -69536

# Row 76: The Monuments Mages
$ws.Range("H76").Value = 69000
$ws.Range("J76").Value = 69000
$ws.Range("L76").Value = 69000
$ws.Range("N76").Value = -69630

# Row 79: Deal with It (L)
$ws.Range("H79").Value = 69000
$ws.Range("J79").Value = 69000
$ws.Range("L79").Value = 69000
$ws.Range("N79").Value = -71184

# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 3091.5173
$ws.Range("J80").Value = 3221.3845
$ws.Range("L80").Value = 3221.3845
$ws.Range("N80").Value = -5217.3845

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 3091.5173
$ws.Range("J83").Value = 3221.3845
$ws.Range("L83").Value = 16106.9225
$ws.Range("N83").Value = -26090.9225

# Row 123: Workplace Workout
$ws.Range("H123").Value = 38898.09
$ws.Range("J123").Value = 38898.09
$ws.Range("L123").Value = 38898.09
$ws.Range("N123").Value = -43798.09

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 2434.2856
$ws.Range("J126").Value = 3562.5
$ws.Range("L126").Value = 10687.5
$ws.Range("N126").Value = -15627.5

# Row 132: On Board for Lar
$ws.Range("H132").Value = 3840.6365
$ws.Range("I132").Value = 2624.6667
$ws.Range("J132").Value = 5299.8
$ws.Range("K132").Value = 7874.000100000001
$ws.Range("L132").Value = 15899.4
$ws.Range("M132").Value = -5344.000100000001
$ws.Range("N132").Value = -20959.4

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")

# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 4285.6763
$ws.Range("I82").Value = 3973.5
$ws.Range("J82").Value = 4731.643
$ws.Range("K82").Value = 3973.5
$ws.Range("L82").Value = 4731.643
$ws.Range("M82").Value = -3612.5
$ws.Range("N82").Value = -5453.643

# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 4285.6763
$ws.Range("I85").Value = 3973.5
$ws.Range("J85").Value = 4731.643
$ws.Range("K85").Value = 3973.5
$ws.Range("L85").Value = 4731.643
$ws.Range("M85").Value = -2725.5
$ws.Range("N85").Value = -7227.643

# Row 100: Tiger in the Sack
$ws.Range("H100").Value = 4550863.5
$ws.Range("I100").Value = 8337833
$ws.Range("K100").Value = 8337833
$ws.Range("M100").Value = -8337292

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 10214.368
$ws.Range("I132").Value = 10214.368
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 30643.104
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -28113.104
$ws.Range("N132").ClearContents()

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")

# Row 51: After the Smock-down
$ws.Range("H51").Value = 25995
$ws.Range("J51").Value = 25995
$ws.Range("L51").Value = 25995
$ws.Range("N51").Value = -27015

# Row 54: No Country for Cold Men
$ws.Range("H54").Value = 28000
$ws.Range("J54").Value = 28000
$ws.Range("L54").Value = 28000
$ws.Range("N54").Value = -29040

# Row 62: Pride Up in Smoke
$ws.Range("H62").Value = 6500
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 6500
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 6500
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -7748

# Row 65: Desperate for Diversionaries (L)
$ws.Range("H65").Value = 6500
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 6500
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 32500
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -38740

# Row 70: An Account of My Boots
$ws.Range("H70").Value = 29995
$ws.Range("I70").Value = 20000
$ws.Range("J70").Value = 39990
$ws.Range("K70").Value = 20000
$ws.Range("L70").Value = 39990
$ws.Range("M70").Value = -19685
$ws.Range("N70").Value = -40620

# Row 73: Soot in My Hair and Scars on My Feet (L)
$ws.Range("H73").Value = 29995
$ws.Range("I73").Value = 20000
$ws.Range("J73").Value = 39990
$ws.Range("K73").Value = 20000
$ws.Range("L73").Value = 39990
$ws.Range("M73").Value = -18908
$ws.Range("N73").Value = -42174

# Row 107: Flax Wax
$ws.Range("H107").Value = 668.0357
$ws.Range("I107").Value = 658.8946999999999
$ws.Range("J107").Value = 687.3333
$ws.Range("K107").Value = 1976.6841
$ws.Range("L107").Value = 2061.9999
$ws.Range("M107").Value = -56.68409999999994
$ws.Range("N107").Value = -5901.9999

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 3937.6287
$ws.Range("I132").Value = 2935.88
$ws.Range("J132").Value = 6442
$ws.Range("K132").Value = 8807.639999999999
$ws.Range("L132").Value = 19326
$ws.Range("M132").Value = -6277.639999999999
$ws.Range("N132").Value = -24386

